# SizingCalculator.xlsx update
#  - bump Django replica count (B19) from 1 to 2 on both sheets
#  - rename the "*_MEMORY_MAX" generated env-var labels to "*_MEMORY_MIN"
#    (they already reference the "Request"/min column D, the label was wrong)
#  - leave the workbook scrolled/selected on the output block (A30:A59),
#    with "UAT & PROD" as the active tab

$wb = $excel.ActiveWorkbook

$sheets = @($wb.Worksheets.Item("DEV"), $wb.Worksheets.Item("UAT & PROD"))

foreach ($ws in $sheets) {
    # Django replicas (min) input: 1 -> 2. Everything downstream (G19, H19,
    # L19, O19, A30, ...) is formula-driven and recalculates automatically.
    $ws.Range("B19").Value = 2

    # The four "_MEMORY_MAX=" labels actually surface column D (the
    # request/min values), so relabel them "_MEMORY_MIN=" to match.
    $ws.Range("A34").Formula = '="DJANGO_MEMORY_MIN="&D12'
    $ws.Range("A40").Formula = '="CELERY_MEMORY_MIN="&D13'
    $ws.Range("A46").Formula = '="REDIS_MEMORY_MIN="&D14'
    $ws.Range("A52").Formula = '="DJANGODB_MEMORY_MIN="&D15'
    $ws.Range("A58").Formula = '="VECTORDB_MEMORY_MIN="&D16'
}

# Restore the selection/active-tab state: both sheets end up scrolled to the
# generated env-var block, with "UAT & PROD" left as the active sheet/tab.
$wsDev = $wb.Worksheets.Item("DEV")
[void]$wsDev.Activate()
$wsDev.Range("A30:A59").Select() | Out-Null

$wsUat = $wb.Worksheets.Item("UAT & PROD")
[void]$wsUat.Activate()
$wsUat.Range("A30:A59").Select() | Out-Null
